$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A and append the new timestamp entry there
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-07-08 22:15:29"
